$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the title text (remove trailing newline) - this causes the shared string
# to be re-inserted at the end of the table when Excel rewrites it.
$title = "8.3.1.2 Экономикадагы иш менен камсыз болгон бардык калктын чакан жана орто ишканаларда иштегендердин үлүшү"
$ws.Range("A1").Value = $title

# Row 1 height and span changes (47.25 -> 45)
$ws.Rows.Item(1).RowHeight = 45

# Row 5 / Row 6 height changes (18.75 -> 17.25)
$ws.Rows.Item(5).RowHeight = 17.25
$ws.Rows.Item(6).RowHeight = 17.25

# Add new column N data (year 2023)
$ws.Cells.Item(4, 14).Value = 2023
$ws.Cells.Item(5, 14).Value = 2.5449890821474286
$ws.Cells.Item(6, 14).Value = 1.4569686017619159

# Copy styles from column M to column N for rows 3-6
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("M5").Copy() | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null
$ws.Range("M6").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null

# Remove selection pane references by resetting selection to A1
$ws.Range("A1").Select() | Out-Null
